$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price quotes in column D (text-typed cells).
# NumberFormat "@" forces the numeric-looking strings to stay text,
# matching how the rest of the Price column is stored.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.24"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.485"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06271"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.653"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.669"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.398"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8335"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01384"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1622"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08298"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03447"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03122"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09301"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.854"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001653"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04778"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006349"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005679"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.713"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3344"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04715"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007043"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1160"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01214"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006263"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004147"

Write-Output "Updated 26 price cells"
